$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the data block (row 43),
# pushing the existing rows 43-69 down to 44-70.
$ws.Rows(43).Insert()

$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44438
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100112012
$ws.Range("G43").Value = "Espinaca"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 30
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 9000
$ws.Range("M43").Value = 9000
$ws.Range("N43").Value = "`$/docena de atados"
$ws.Range("O43").Value = "Región de La Araucanía"
$ws.Range("P43").Value = 3000
$ws.Range("Q43").Value = 3
$ws.Range("R43").Value = "Hortaliza"

$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat
